# Route Mobile Ltd - "Quarterly" sheet: split the combined "Mon YY Qn" label
# column into three separate columns (Year / Month / Quarter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# ---------------------------------------------------------------------------
# 1. Insert two new columns right after column A for the new "Month" and
#    "Quarter" fields (existing columns B..U shift to D..W).
# ---------------------------------------------------------------------------
$ws.Range("B1:C1").EntireColumn.Insert()

# Header labels for the two new columns.
$ws.Range("B1").Value = "Month"
$ws.Range("C1").Value = "Quarter"

# ---------------------------------------------------------------------------
# 2. Per-row Year / Month / Quarter values, derived from the old combined
#    "Mon YY Qn" labels that used to live in column A.
# ---------------------------------------------------------------------------
$data = @(
  @{Row=2; Year="2019"; Month="12"; Quarter="Q3"},
  @{Row=3; Year="2020"; Month="03"; Quarter="Q4"},
  @{Row=4; Year="2020"; Month="06"; Quarter="Q1"},
  @{Row=5; Year="2020"; Month="09"; Quarter="Q2"},
  @{Row=6; Year="2020"; Month="12"; Quarter="Q3"},
  @{Row=7; Year="2021"; Month="03"; Quarter="Q4"},
  @{Row=8; Year="2021"; Month="06"; Quarter="Q1"},
  @{Row=9; Year="2021"; Month="09"; Quarter="Q2"},
  @{Row=10; Year="2021"; Month="12"; Quarter="Q3"},
  @{Row=11; Year="2022"; Month="03"; Quarter="Q4"},
  @{Row=12; Year="2022"; Month="06"; Quarter="Q1"},
  @{Row=13; Year="2022"; Month="09"; Quarter="Q2"},
  @{Row=14; Year="2022"; Month="12"; Quarter="Q3"},
  @{Row=15; Year="2023"; Month="03"; Quarter="Q4"},
  @{Row=16; Year="2023"; Month="06"; Quarter="Q1"},
  @{Row=17; Year="2023"; Month="09"; Quarter="Q2"},
  @{Row=18; Year="2023"; Month="12"; Quarter="Q3"},
  @{Row=19; Year="2024"; Month="03"; Quarter="Q4"},
  @{Row=20; Year="2024"; Month="06"; Quarter="Q1"},
  @{Row=21; Year="2024"; Month="09"; Quarter="Q2"}
)

# Year (col A) and Month (col B) are numeric-looking strings ("2019", "03", ...).
# Writing them straight to .Value would auto-coerce to a number (and drop the
# leading zero on the month code), so stage them as TEXT() formulas first and
# immediately flatten the range to literal values via copy / paste-special.
# This keeps the cells as plain shared-string text without picking up a
# "quote-prefixed" style (which a leading apostrophe would add).
foreach ($item in $data) {
  $r = $item.Row
  $ws.Cells.Item($r, 1).Formula = '=TEXT(' + $item.Year + ',"0")'
  $ws.Cells.Item($r, 2).Formula = '=TEXT(' + $item.Month + ',"00")'
}

$yearMonthRange = $ws.Range("A2:B21")
$yearMonthRange.Copy()
$yearMonthRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# Quarter (col C) is plain alphanumeric text ("Q1".."Q4") - no numeric
# coercion risk, so it can be written directly.
foreach ($item in $data) {
  $ws.Cells.Item($item.Row, 3).Value = $item.Quarter
}
